$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Checkmark character used throughout the sheet (Wingdings-rendered U+00FC)
$check = [string][char]0x00FC

# Row 9 / 10 — new "GSB Frais" / "GSB Visites" lines with checkmarks
$ws.Range("A9").Value = "GSB Frais"
$ws.Range("F9").Value = $check
$ws.Range("G9").Value = $check

$ws.Range("A10").Value = "GSB Visites"
$ws.Range("F10").Value = $check

# Row 21 / 22 — new realisation lines with checkmarks
$ws.Range("A21").Value = "Mise en place d" + [string][char]0x2019 + "un trello et UML + user stories"
$ws.Range("F21").Value = $check

$ws.Range("A22").Value = "Mise en place de WordPress non-installé sur le PC"
$ws.Range("H22").Value = $check

# Row 28 — existing text amended
$ws.Range("A28").Value = "J'ai amélioré le site internet interne du Médipôle de Savoie (trello, user stories)"

# Restore the selection recorded in the saved workbook
$ws.Range("J24").Select() | Out-Null
